# "Added Emailed Reciepts stretch goal."
# Adds a new "Expense" worksheet (after "Sarah") that itemises each
# client's order into client/product/price/subtotal rows, and updates
# the active-sheet/selection state to match.

$wb = $excel.ActiveWorkbook
$sarah = $wb.Worksheets.Item("Sarah")

# Add the new "Expense" sheet directly after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$expense = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$expense.Name = "Expense"

# Write the expense line items first (so new shared strings are interned
# in client/product order), then the header row.
$rows = @(
    @(1, "Emma",  "Emma_email",  13,  "Store",                      3.5, 7),
    @(1, "Emma",  "Emma_email",  13,  "Taco",                       2,   6),
    @(2, "Tammy", "Tammy_email", 32,  "Tuba",                       4,   32),
    @(3, "Sarah", "Sarah_Email", 281, "Chai",                       2,   8),
    @(3, "Sarah", "Sarah_Email", 281, "Northwoods Cranberry Sauce", 40,  240),
    @(3, "Sarah", "Sarah_Email", 281, "Store",                      3.5, 21),
    @(3, "Sarah", "Sarah_Email", 281, "Tuba",                       4,   12)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Count; $c++) {
        $expense.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

$headers = @("clientid", "clientname", "clientemail", "total", "productname", "price", "subtotal")
for ($c = 1; $c -le $headers.Count; $c++) {
    $expense.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Restore per-sheet selections; activating Expense last makes it the
# selected/visible tab, matching the saved workbook state.
$sarah.Range("C14").Select()
$expense.Range("C38").Select()
